# Insert a new data row at row 636, shifting existing rows 636-672 down to 637-673.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before row 636 (this pushes old row 636.. down to 637..)
$ws.Rows.Item(636).Insert()

# Populate the new row 636 with the new weekly data entry.
$ws.Cells.Item(636, 1).Value = 6
$ws.Cells.Item(636, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(636, 3).Value = "Metropolitana"
$ws.Cells.Item(636, 4).Value = 45132
$ws.Cells.Item(636, 5).Value = 13
$ws.Cells.Item(636, 6).Value = "Fruta"
$ws.Cells.Item(636, 7).Value = 100107
$ws.Cells.Item(636, 8).Value = "Otros"
$ws.Cells.Item(636, 9).Value = 100107011
$ws.Cells.Item(636, 10).Value = "Tuna"
$ws.Cells.Item(636, 11).Value = "Sin especificar"
$ws.Cells.Item(636, 12).Value = "Especial"
$ws.Cells.Item(636, 13).Value = 120
$ws.Cells.Item(636, 14).Value = 30000
$ws.Cells.Item(636, 15).Value = 30000
$ws.Cells.Item(636, 16).Value = 30000
$ws.Cells.Item(636, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(636, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(636, 19).Value = 1667
$ws.Cells.Item(636, 20).Value = 18
